$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2117647058823529
$ws.Range("C2").Value = 0.5058823529411764
$ws.Range("J2").Value = 0.02941176470588235
$ws.Range("P2").Value = 0.1235294117647059
$ws.Range("S2").Value = 0.1294117647058824
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.01675977653631285
$ws.Range("J3").Value = 0.03910614525139665
$ws.Range("P3").Value = 0.7262569832402235
$ws.Range("S3").Value = 0.2122905027932961
$ws.Range("J4").Value = 0.06976744186046512
$ws.Range("P4").Value = 0.5813953488372093
$ws.Range("S4").Value = 0.3488372093023256
$ws.Range("B6").Value = 0.06167400881057269
$ws.Range("D6").Value = 0.00881057268722467
$ws.Range("F6").Value = 0.04845814977973568
$ws.Range("J6").Value = 0.3656387665198238
$ws.Range("O6").Value = 0.03524229074889868
$ws.Range("Q6").Value = 0.13215859030837
$ws.Range("R6").Value = 0.04845814977973568
$ws.Range("S6").Value = 0.2995594713656388
$ws.Range("B7").Value = 0.1019607843137255
$ws.Range("D7").Value = 0.00392156862745098
$ws.Range("E7").Value = 0.007843137254901961
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1490196078431373
$ws.Range("O7").Value = 0.0196078431372549
$ws.Range("Q7").Value = 0.1372549019607843
$ws.Range("R7").Value = 0.08235294117647059
$ws.Range("S7").Value = 0.4392156862745098
$ws.Range("B8").Value = 0.1002044989775051
$ws.Range("D8").Value = 0.016359918200409
$ws.Range("F8").Value = 0.06339468302658487
$ws.Range("J8").Value = 0.114519427402863
$ws.Range("O8").Value = 0.0245398773006135
$ws.Range("Q8").Value = 0.16359918200409
$ws.Range("R8").Value = 0.0736196319018405
$ws.Range("S8").Value = 0.4437627811860941
$ws.Range("B9").Value = 0.125
$ws.Range("D9").Value = 0.02083333333333333
$ws.Range("J9").Value = 0.15625
$ws.Range("O9").Value = 0.005208333333333333
$ws.Range("Q9").Value = 0.1354166666666667
$ws.Range("R9").Value = 0.07291666666666667
$ws.Range("S9").Value = 0.4010416666666667
$ws.Range("B10").Value = 0.1102719033232628
$ws.Range("D10").Value = 0.02341389728096677
$ws.Range("F10").Value = 0.05211480362537765
$ws.Range("J10").Value = 0.1072507552870091
$ws.Range("O10").Value = 0.01661631419939577
$ws.Range("Q10").Value = 0.2401812688821752
$ws.Range("R10").Value = 0.07854984894259819
$ws.Range("S10").Value = 0.3716012084592145
$ws.Range("G11").Value = 0.1862244897959184
$ws.Range("J11").Value = 0.07653061224489796
$ws.Range("K11").Value = 0.2270408163265306
$ws.Range("L11").Value = 0.4872448979591837
$ws.Range("S11").Value = 0.02295918367346939
$ws.Range("G12").Value = 0.7437185929648241
$ws.Range("J12").Value = 0.1909547738693467
$ws.Range("L12").Value = 0.02512562814070352
$ws.Range("S12").Value = 0.04020100502512563
$ws.Range("G13").Value = 0.6779661016949152
$ws.Range("J13").Value = 0.3220338983050847
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1650485436893204
$ws.Range("I15").Value = 0.04368932038834952
$ws.Range("J15").Value = 0.354368932038835
$ws.Range("K15").Value = 0.07281553398058252
$ws.Range("M15").Value = 0.004854368932038835
$ws.Range("O15").Value = 0.01456310679611651
$ws.Range("S15").Value = 0.325242718446602
$ws.Range("F16").Value = 0.0154639175257732
$ws.Range("H16").Value = 0.1752577319587629
$ws.Range("I16").Value = 0.08247422680412371
$ws.Range("J16").Value = 0.3969072164948453
$ws.Range("K16").Value = 0.134020618556701
$ws.Range("M16").Value = 0.01030927835051546
$ws.Range("O16").Value = 0.04123711340206185
$ws.Range("S16").Value = 0.1443298969072165
$ws.Range("F17").Value = 0.02484472049689441
$ws.Range("H17").Value = 0.1635610766045549
$ws.Range("I17").Value = 0.09109730848861283
$ws.Range("J17").Value = 0.401656314699793
$ws.Range("K17").Value = 0.113871635610766
$ws.Range("M17").Value = 0.02484472049689441
$ws.Range("O17").Value = 0.03519668737060042
$ws.Range("S17").Value = 0.1449275362318841
$ws.Range("F18").Value = 0.04324324324324325
$ws.Range("H18").Value = 0.1675675675675676
$ws.Range("I18").Value = 0.06486486486486487
$ws.Range("J18").Value = 0.3675675675675676
$ws.Range("K18").Value = 0.1351351351351351
$ws.Range("M18").Value = 0.01081081081081081
$ws.Range("O18").Value = 0.06486486486486487
$ws.Range("S18").Value = 0.145945945945946
$ws.Range("F19").Value = 0.02356202356202356
$ws.Range("H19").Value = 0.2210672210672211
$ws.Range("I19").Value = 0.07900207900207901
$ws.Range("J19").Value = 0.3367983367983368
$ws.Range("K19").Value = 0.1191961191961192
$ws.Range("M19").Value = 0.0297990297990298
$ws.Range("N19").Value = 0.001386001386001386
$ws.Range("O19").Value = 0.05890505890505891
$ws.Range("S19").Value = 0.1302841302841303
